# "Generate Report for Handoff"
# A new handoff run replaces the old report GUID/content-hash pair with a
# freshly generated one and bumps the recorded handoff timestamps.

$wb = $excel.ActiveWorkbook

$newGuid = "b28b68be-e023-406f-9538-00996f750331"
$newHash = "bf4c7b0d12fb7a6574002a2398cc160ad8194340"

$newMd    = "$newGuid.md"
$newZhXlf = "$newGuid.$newHash.zh-cn.xlf"
$newDeXlf = "$newGuid.$newHash.de-de.xlf"

# ---------------------------------------------------------------------
# Overview sheet: handoff file name + latest handoff date
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newMd
$wsOverview.Range("D2").Value = "2016-49-21 02:49:19"
# keep the hyperlink's visible text in sync with the new file name
$wsOverview.Hyperlinks.Item(1).TextToDisplay = $newMd

# ---------------------------------------------------------------------
# zh-cn sheet: handoff file name, target xlf name + its handoff datetime
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = $newMd
$wsZh.Range("D2").Value = $newZhXlf
$wsZh.Range("E2").Value = "2016-03-21 02:49:16"
# Hyperlinks, in sheet order: A2 (1), B2 (2, unchanged ".md"), D2 (3)
$wsZh.Hyperlinks.Item(1).TextToDisplay = $newMd
$wsZh.Hyperlinks.Item(3).TextToDisplay = $newZhXlf

# ---------------------------------------------------------------------
# de-de sheet: handoff file name, target xlf name + its handoff datetime
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = $newMd
$wsDe.Range("D2").Value = $newDeXlf
$wsDe.Range("E2").Value = "2016-03-21 02:49:19"
# Hyperlinks, in sheet order: A2 (1), B2 (2, unchanged ".md"), D2 (3)
$wsDe.Hyperlinks.Item(1).TextToDisplay = $newMd
$wsDe.Hyperlinks.Item(3).TextToDisplay = $newDeXlf
